$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.4304569862216283, 0.1468426154918987, 0.05358668348697648, 0.1206899266050296, 1.019680491591984, 0.9743034593364683, 0.2759078690222339, 0.2112346193438484, 0.1391969961823882, 3.673779994777647),
    @(0.3949936695926226, 0.1450407181110336, 0.05133406264843643, 0.1212344157328804, 1.022163523110493, 0.9820457855585438, 0.2421641450224001, 0.2086123896109626, 0.1323664944757788, 3.695655788436966),
    @(0.3733011042209569, 0.1439243604565306, 0.04993567548891775, 0.1216141015284951, 1.024239218371036, 0.9872413179532025, 0.2214114467003583, 0.2071037662165907, 0.1282237426441526, 3.710914529281411),
    @(0.3644824466778402, 0.1434669499554531, 0.04936200680958081, 0.1217802536036832, 1.025223740669638, 0.9894696601875701, 0.2129465139258571, 0.2065145667823813, 0.1265485363598629, 3.717592163081008),
    @(0.363019415951527, 0.1433908479062325, 0.04926651990710695, 0.1218085337538621, 1.025395596899735, 0.9898463875514878, 0.2115404471205835, 0.2064182775000489, 0.1262711583505869, 3.718728741477122),
    @(0.3731820859210586, 0.1439182016893312, 0.04992795420470486, 0.1216162960204912, 1.024251934449985, 0.9872709201692302, 0.2212973175604134, 0.2070957164143294, 0.1282010974314893, 3.711002725156021),
    @(0.4182125558631355, 0.1462234083143983, 0.05281316320574092, 0.1208682618361259, 1.020422305380521, 0.9768813649246546, 0.2642803546536925, 0.2103094620510504, 0.136831283284323, 3.680943697643599),
    @(0.5071461540983648, 0.150663775998126, 0.05834904085697445, 0.1197605389680216, 1.017282644966187, 0.9600102603098222, 0.3482834488540334, 0.2174139395231549, 0.1541568723306952, 3.636487554016043),
    @(0.5728461095904152, 0.1538763405016255, 0.06234106942189754, 0.1191646370865609, 1.017638205870142, 0.9497479141043144, 0.409808439535766, 0.2231201489645542, 0.1671263506892871, 3.612652957088358),
    @(0.6028085274433579, 0.15532683600631, 0.06414067257035327, 0.1189406587720594, 1.018377728254485, 0.9455419240679745, 0.4377527216997237, 0.225821204556695, 0.1730778080252406, 3.603725710979546),
    @(0.6141648029473004, 0.1558745098908645, 0.06481975759317748, 0.1188625988993675, 1.018740795757566, 0.9440156794015451, 0.4483277895185722, 0.2268591028726803, 0.1753387809628606, 3.6006204869974),
    @(0.6117185844158541, 0.1557566299975406, 0.06467361104071045, 0.1188791102939799, 1.018658910971929, 0.9443414274435611, 0.4460505715160821, 0.2266349039127107, 0.1748515180671575, 3.601277008399393),
    @(0.6037426152118144, 0.1553719256346824, 0.06419658927983107, 0.1189341014401659, 1.018405934393734, 0.9454150269106023, 0.4386228787441269, 0.2259062915132404, 0.1732636743757325, 3.6034647247522),
    @(0.5988584073218419, 0.1551360742764842, 0.06390408814583282, 0.1189686644082322, 1.018261789715169, 0.946081293654224, 0.4340722977316318, 0.2254619554669119, 0.1722920191162771, 3.60484061800318),
    @(0.5708894520985552, 0.1537813251896338, 0.06222312884982983, 0.1191802209127548, 1.017601497819477, 0.950032088420766, 0.407981288966738, 0.2229457409250841, 0.1667384355704513, 3.613274906369952),
    @(0.5537501943411201, 0.1529474149506882, 0.06118769631697063, 0.1193220553358163, 1.01734437215147, 0.95257419092313, 0.391963715106499, 0.2214290372116636, 0.163344608392471, 3.618939551067854),
    @(0.5438992616703615, 0.1524667459892868, 0.06059060185559417, 0.1194080690850932, 1.01725086674999, 0.9540798658029033, 0.3827467447112554, 0.220566582324949, 0.1613974312466553, 3.622377980246426),
    @(0.5405651470556734, 0.1523038243892643, 0.06038817212742487, 0.1194379539869388, 1.01722855008083, 0.9545971374896638, 0.3796253515259878, 0.2202762746381808, 0.1607389896648712, 3.623573138103922),
    @(0.5555739659360768, 0.153036292490988, 0.06129807957891131, 0.119306498006031, 1.017366115230395, 0.9522990752704459, 0.3936692399072967, 0.2215894674688599, 0.16370538479309, 3.618317882699387),
    @(0.6060850802116988, 0.1554849661990403, 0.06433676710750547, 0.1189177659920286, 1.018477986893615, 0.9450978808421127, 0.4408047616568638, 0.2261198942834, 0.1737298658864503, 3.602814667666365),
    @(0.639155920891028, 0.1570759931764627, 0.06630880648236115, 0.1187030774923521, 1.01968857362219, 0.9407789269421798, 0.4715705190787958, 0.2291685801826446, 0.1803238466073509, 3.594287219584231),
    @(0.6215002160052165, 0.1562276942347722, 0.06525757497593787, 0.1188140640540585, 1.018998201635185, 0.9430485896413145, 0.4551540945910233, 0.2275334297954856, 0.1768006785427048, 3.598691665875492),
    @(0.5547494304794895, 0.1529961148258252, 0.06124818094355788, 0.1193135175462494, 1.017356115978536, 0.9524233174174057, 0.3928981978148158, 0.2215169072963192, 0.163542265377707, 3.618598372822476),
    @(0.4830223065539201, 0.1494712137925021, 0.05686458403999239, 0.1200218609907768, 1.017664340423231, 0.9641996066271119, 0.3255908650441199, 0.2154063605009782, 0.1494272990656249, 3.646963519823402)
)

$colLetters = @("B","C","D","E","F","I","K","L","M","O")
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    for ($c = 0; $c -lt $colLetters.Length; $c++) {
        $ws.Range($colLetters[$c] + $rowNum).Value = $data[$r][$c]
    }
}
